$d = $word.ActiveDocument

# 1. Split "Лабораторна робота №3" into "Лабораторна робота №" + "4" (two runs)
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Лабораторна робота №3", $true, $false, $false, $false, $false, $true, 1, $false, "Лабораторна робота №4", 2)

# 2. Header "Лабораторна робота №3" -> "...4" (page? it's in header1.xml)
$headerFind = $word.ActiveDocument.Sections(1).Headers(1).Range.Find
$headerFind.ClearFormatting()
$headerFind.Execute("робота №3", $true, $false, $false, $false, $false, $true, 1, $false, "робота №4", 2)
